# 13. Hafta Butce Guncellemesi
# Adds the week-13 ("Son Teslim") budget row + submit-log row, rolls the
# running totals forward one row, bumps the submit date to 22 May 2015,
# and drops a footnote marker "5" next to the new entry — mirroring the
# markers already sitting next to weeks 2/4/6/9.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Submit date (top of sheet) moves from 15 May 2015 to 22 May 2015.
# ---------------------------------------------------------------------
$ws.Range("N9").Formula = "=DATE(2015,5,22)"

# ---------------------------------------------------------------------
# 2) Fill in week 13's budget figures (row 44). The row already carried
#    its label ("15.05.15 - 22.05.15") and formatting; only the numbers
#    were missing.
# ---------------------------------------------------------------------
$ws.Range("C44").Value = 10000
$ws.Range("D44").Value = 6000
$ws.Range("E44").Value = 6000
$ws.Range("F44").Value = 4000
$ws.Range("G44").Value = 4000
$ws.Range("H44").Value = 4000
$ws.Range("I44").Value = 0
$ws.Range("J44").Formula = "=SUM(C44:I44)"
$ws.Range("K44").Value = 100000
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = 0
$ws.Range("N44").Formula = "=( (M44 * L44) * K44 / 100 ) + K44"
$ws.Range("O44").Formula = "=N44 - J44"
$ws.Range("P44").Formula = "=Q44 * 0.1 * 0"
$ws.Range("Q44").Formula = "= (Q43 + O44) + P43"

# ---------------------------------------------------------------------
# 3) "Toplam" row (45) now needs to roll up through row 44 instead of 43.
# ---------------------------------------------------------------------
$ws.Range("J45").Formula = "=SUM(J32:J44)"
$ws.Range("N45").Formula = "=SUM(N32:N44)"
$ws.Range("P45").Formula = "=SUM(P32:P44) - P44"
$ws.Range("Q45").Formula = "= Q44"

# ---------------------------------------------------------------------
# 4) Extend the "SUBMIT BILGILENDIRME" table with the week-13 submission
#    (row 57): NO=5, TARIH=22.05.2015, HAFTA=13, ACIKLAMA="Son Teslim".
# ---------------------------------------------------------------------
$tbl = $ws.ListObjects.Item(1)
$tbl.ListRows.Add() | Out-Null

# carry the same formatting the previous row (56) used
$ws.Range("G56:J56").Copy()
$ws.Range("G57:J57").PasteSpecial(-4122)

$ws.Range("G57").Value = 5
$ws.Range("H57").Formula = "=DATE(2015,5,22)"
$ws.Range("I57").Value = 13
$ws.Range("J57").Value = "Son Teslim"

# ---------------------------------------------------------------------
# 5) Drop a new footnote textbox ("5") next to the new row, matching the
#    placement/size of the existing 1-4 markers beside weeks 2/4/6/9.
# ---------------------------------------------------------------------
$shp = $ws.Shapes.AddTextbox(1, 941.2585611466535, 645.9881102362205, 12.405354330708661, 19.59456692913386)
$shp.TextFrame.Characters().Text = "5"
$shp.TextFrame.Characters().Font.Italic = $true
$shp.TextFrame.Characters().Font.Size = 10

# ---------------------------------------------------------------------
# 6) Move the saved selection to reflect where the user ended up editing.
# ---------------------------------------------------------------------
$ws.Range("J58").Select()
